# Port battle calculator correction
# - Reshuffle the ship list on the "Shallow water port" sheet (rows 4-21)
#   so that the ship names line up with their correct BR values again.
# - Expand the BR-total SUM() ranges on both sheets to include the rows
#   that were already present (row 35 on "Deep water port", row 21 on
#   "Shallow water port") but had been left out of the totals.

$wb = $excel.ActiveWorkbook

$wsDeep = $wb.Worksheets.Item("Deep water port")
$wsShallow = $wb.Worksheets.Item("Shallow water port")

# --- Shallow water port: fix ship / BR alignment for rows 4-21 ---------
$ships = @(
    @{ Row = 4;  Name = "Hercules";             BR = 100 },
    @{ Row = 5;  Name = "Pandora";               BR = 100 },
    @{ Row = 6;  Name = "Mercury";               BR = 80 },
    @{ Row = 7;  Name = "Mortar Brig";           BR = 80 },
    @{ Row = 8;  Name = "NavyBrig";              BR = 80 },
    @{ Row = 9;  Name = "Niagara";               BR = 80 },
    @{ Row = 10; Name = "Prince de Neufchatel";  BR = 80 },
    @{ Row = 11; Name = "Rattlesnake";           BR = 80 },
    @{ Row = 12; Name = "Rattlesnake Heavy";     BR = 80 },
    @{ Row = 13; Name = "Snow";                  BR = 80 },
    @{ Row = 14; Name = "Brig";                  BR = 70 },
    @{ Row = 15; Name = "Pickle";                BR = 55 },
    @{ Row = 16; Name = "Cutter";                BR = 50 },
    @{ Row = 17; Name = "GunBoat";               BR = 50 },
    @{ Row = 18; Name = "Lynx";                  BR = 50 },
    @{ Row = 19; Name = "Privateer";             BR = 50 },
    @{ Row = 20; Name = "Yacht";                 BR = 50 },
    @{ Row = 21; Name = "Yacht Silver";          BR = 50 }
)

foreach ($ship in $ships) {
    $wsShallow.Cells.Item($ship.Row, 2).Value = $ship.Name
    $wsShallow.Cells.Item($ship.Row, 3).Value = $ship.BR
}

# --- Expand BR-total formulas to include the already-present last row --
$wsDeep.Range("D3").Formula = "=SUM(D4:D35)"
$wsDeep.Range("E3").Formula = "=SUM(E4:E35)"

$wsShallow.Range("D3").Formula = "=SUM(D4:D21)"
$wsShallow.Range("E3").Formula = "=SUM(E4:E21)"
